# Updated cryptos list with GitHub Actions
# Each text cell is written with a leading literal apostrophe so Excel keeps
# it as text (it would otherwise auto-convert numeric-looking strings like
# "211.05" or "22.70" into numbers, losing exact formatting such as trailing
# zeros). The Style is then reset to "Normal" so no stray cell formatting
# (e.g. quote-prefix) is left applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.510.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.78%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.619.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.73%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.37%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''211.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.91%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.523'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.60%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.26%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''22.70'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -1.92%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.263'
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = '''  +0.02%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0885'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -0.47%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''1.847.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -1.51%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.618.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.43%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  -0.55%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  -2.16%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '''  +1.26%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''27.518.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -0.40%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  -0.27%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  -0.67%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  -2.21%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  +0.16%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  -0.59%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''10.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.67%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  +6.44%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''149.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.40%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -1.72%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  +0.30%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''6.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.46%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''15.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -0.45%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -0.55%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.0482'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -0.82%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -1.29%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''1.446.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +0.76%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''3.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -3.68%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -4.03%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.12%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.938'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +4.59%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.561'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -2.38%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  -0.06%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.862'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -2.39%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''69.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +5.74%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  +0.25%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -3.06%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  -0.29%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''5.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.51%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  -2.37%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''1.758.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -1.41%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -0.03%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''86.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -0.17%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -0.88%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.0996'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.64%  '
$ws.Range("E51").Style = "Normal"
